# correção das notas do fórum para matc65 em 2021.2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Zero out all forum-view/score data (columns B..J, rows 2..50),
# leaving the "matricula" column (A) and header row (1) untouched.
$ws.Range("B2:J50").Value = 0
